$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 793 entirely, shifting rows 794-804 up by one.
$ws.Rows.Item(793).Delete()
